$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.119.76'
$ws.Range("E2").Value = '  +1.48%  '
$ws.Range("D3").Value = '3.267.86'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.88'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.25'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.83%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -1.05%  '
$ws.Range("E9").Value = '  +3.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.72'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.417'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("D12").Value = '3.832.88'
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.64'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.58%  '
$ws.Range("D15").Value = '68.134.38'
$ws.Range("E15").Value = '  +1.46%  '
$ws.Range("E16").Value = '  +2.49%  '
$ws.Range("D17").Value = '3.262.45'
$ws.Range("E17").Value = '  -0.26%  '
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("E19").Value = '  +1.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '382.81'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.71'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.33'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.516'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000121'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.88'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  +2.92%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("E29").Value = '  +0.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.75'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.25'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.92%  '
$ws.Range("E32").Value = '  +1.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.28'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("E35").Value = '  +2.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.87'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.47%  '
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.837'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.76'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.68'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.62'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.61'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.74%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.47'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '348.65'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '25.42'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0687'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.21%  '
$ws.Range("D47").Value = '2.647.80'
$ws.Range("E47").Value = '  -3.66%  '
$ws.Range("E48").Value = '  +0.97%  '
$ws.Range("E49").Value = '  +5.00%  '
$ws.Range("E50").Value = '  -0.74%  '
$ws.Range("E51").Value = '  +1.05%  '
